$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.255.47"
$ws.Range("E2").Value = "  +5.40%  "
$ws.Range("D3").Value = "2.745.75"
$ws.Range("E3").Value = "  +3.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.07"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.41"
$ws.Range("E6").Value = "  +9.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.611"
$ws.Range("E8").Value = "  +1.85%  "
$ws.Range("D9").Value = "2.773.83"
$ws.Range("E9").Value = "  +3.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.78"
$ws.Range("E10").Value = "  +2.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.114"
$ws.Range("E11").Value = "  +6.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.395"
$ws.Range("E12").Value = "  +3.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.159"
$ws.Range("E13").Value = "  +2.01%  "
$ws.Range("D14").Value = "3.236.67"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.35"
$ws.Range("E15").Value = "  +4.82%  "
$ws.Range("D16").Value = "63.855.38"
$ws.Range("E16").Value = "  +4.80%  "
$ws.Range("E17").Value = "  +7.65%  "
$ws.Range("D18").Value = "2.766.61"
$ws.Range("E18").Value = "  +3.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.14"
$ws.Range("E19").Value = "  +4.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.96"
$ws.Range("E20").Value = "  +4.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "364.65"
$ws.Range("E22").Value = "  +2.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.536"
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.04"
$ws.Range("E25").Value = "  +4.76%  "
$ws.Range("E26").Value = "  +5.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.65"
$ws.Range("E27").Value = "  +6.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "0.0₃0917"
$ws.Range("E29").Value = "  +12.97%  "
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.23"
$ws.Range("E31").Value = "  +5.54%  "
$ws.Range("E32").Value = "  +20.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "174.91"
$ws.Range("E33").Value = "  +6.77%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.78"
$ws.Range("E34").Value = "  +4.26%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.92"
$ws.Range("E36").Value = "  +7.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.46"
$ws.Range("E37").Value = "  +9.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.83"
$ws.Range("E38").Value = "  +10.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.01"
$ws.Range("E39").Value = "  +10.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.33"
$ws.Range("E40").Value = "  +5.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "342.99"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.50"
$ws.Range("E42").Value = "  +2.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.91"
$ws.Range("E43").Value = "  +12.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.19"
$ws.Range("E44").Value = "  +8.93%  "
$ws.Range("E45").Value = "  +7.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0602"
$ws.Range("E46").Value = "  +6.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.655"
$ws.Range("E47").Value = "  +4.96%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "139.16"
$ws.Range("E48").Value = "  +4.50%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0260"
$ws.Range("E49").Value = "  +4.18%  "
$ws.Range("E50").Value = "  +2.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("E51").Value = "  -0.06%  "
